$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.587.56"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "2.091.49"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0840"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "2.401.56"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.77%  "
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "2.086.99"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "38.492.31"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("E28").Value = "  +5.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +6.59%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "1.539.58"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0927"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.90%  "
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "2.288.49"
$ws.Range("E51").Value = "  +2.08%  "
